$wb = $excel.ActiveWorkbook

# Update the "Priority" column (E) for rows 8-13 on the zh-cn and de-de sheets
# from empty string to "ht" (matching the other handoff-type rows).
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

foreach ($row in 8..13) {
    $zhcn.Range("E$row").Value = "ht"
    $dede.Range("E$row").Value = "ht"
}

# Update the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# for the 3d48c809 file row across all three sheets.
$overview = $wb.Worksheets.Item("Overview")

foreach ($row in 8..13) {
    $overview.Range("G$row").Value = "2016-09-01 02:26:14"
    $dede.Range("H$row").Value = "2016-09-01 02:26:14"
    $zhcn.Range("H$row").Value = "2016-09-01 02:26:02"
}
